$d = $word.ActiveDocument

# --- Simple in-place text replacements (paragraph 1, 2, 4) ---
$d.Paragraphs.Item(1).Range.Text = 'Review 174: In-Context Learning Creates Task Vectors'
$d.Paragraphs.Item(2).Range.Text = 'Paper: https://arxiv.org/abs/2311.06668v3'
$d.Paragraphs.Item(4).Range.Text = 'https://huggingface.co/papers/2310.15916'

# --- Insert 5 new "Normal"-style paragraphs right after paragraph 5 ---
# (paragraph 5 is the existing empty line after the hf link)
for ($i = 0; $i -lt 5; $i++) {
    $d.Paragraphs.Item(5 + $i).Range.InsertParagraphAfter()
}

# Paragraph 6 stays empty (already inserted blank).
$d.Paragraphs.Item(7).Range.Text = 'אחד היכולות המדהימות של מודלי שפה ענקיים היא יכולת למידת in-context או ICL בקצרה. ICL היא יכולת של LLM ללמוד מכמה דוגמאות בלי לשנות בכלל את המשקלים שלו. כלומר אנו מעבירים למודל שפה כמה דוגמאות בסגנון (מלון -> צהוב, מלפפון -> ירוק,..) ולאחר המכן אם תזינו למודל ״בננה -> …״, הוא יבין שמדובר בצבע ויענה צהוב. '
$d.Paragraphs.Item(8).Range.Text = ' '
$d.Paragraphs.Item(9).Range.Text = 'אבל איך המנגנון הזה עובד? המאמר המסוקר טוען ומראה שמדובר כאן בתהליך דו שלבי: '
$d.Paragraphs.Item(10).Range.Text = '– הזנה של הדוגמאות (נסמן אותם ב S) המחשבים את הפרמטרים של פונקציה מסוימת (בהמשך נסביר איך היא בנויה) שתופעל על דוגמת הטסט x (בננה במקרה המתואר). '

# --- Old paragraph 6 ("אוקיי...") is now paragraph 11; replace its text ---
$d.Paragraphs.Item(11).Range.Text = '– הפעלה של פונקציה זו על שאילת טסט x. המאמר טוען שהפרמטרים האלו לא תלויים בשאילתת הטסט x עצמו אלא רק ב- S (במאמר זה מנוסח בצורה מתמטית יפה שמאוד אהבתי). ההשערה הזו היא לא לגמרי טריוויאלית כי בארכיטקטורת הטרנספורמרים הייצוג של דוגמאות מתויגות S תלוי גם בשאילתה x.'
$d.Paragraphs.Item(11).Range.InsertParagraphAfter()

# --- Old paragraph 7 ("גנרוט...") is now paragraph 13; replace its text ---
$d.Paragraphs.Item(13).Range.Text = 'המאמר מראה שב- ICL ניתן להגיע להפרדה כזו בין ייצוג המשימה (הנגזר מ- S) וייצוג השאילתה x. אוקיי, אז מה זה הפרמטרים האלו שמחושבים רק על דוגמאות S? המאמר טוען הם בעצם הפלטים של שכבה L של הטרנספורמר עבור הטוקנים של S כאשר L אינה השכבה האחרונה של מודל השפה. פרמטרים אלו מגדירים(דרך הזנה) לפונקציה שהיא הפעלה של השכבות הנותרות על פלט זה (= ייצוג המשימה) וגם על השאילתה x. '
$d.Paragraphs.Item(13).Range.InsertParagraphAfter()

# --- Old paragraph 8 ("כיול...") is now paragraph 15; replace its text ---
$d.Paragraphs.Item(15).Range.Text = 'איך הם בדקו זאת? אוקיי, השאילתה מורכבת מגוף השאלה (בננה בדוגמה שלנו) ובסימן שאלה מאולתר (״->״ במקרה) שלנו המאותת למודל שפה שהוא צריך לפתור אותה. אז המחברים העתיקו את ייצוג של ״->״ בשכבה L עבור דוגמא לא קשורה x'' ואז ממשיכים עם השאילתה המקורית לאחר מכן.  המאמר מראה שעבור שכבה מסוימת L החלפה כזו לא מובילה לירידה ניכרת בביצועים(יחסית לייצוג של ״->״ הנבנה באופן רגיל). '
$d.Paragraphs.Item(15).Range.InsertParagraphAfter()

# --- Old paragraph 9 ("מפעילים...") is now paragraph 17; delete it entirely ---
$d.Paragraphs.Item(17).Range.Delete()

# --- Old paragraph 10 ("זהו זה...") is now paragraph 17; replace its text ---
$d.Paragraphs.Item(17).Range.Text = 'כלומר הפלט של שכבה L של מודל שפה עבור הטוקנים של S אכן לא תלויה בשאילתה x. מה שמעניין שעבור מודלי שפה בגדלים שונים L האופטימלי יצא בערך 15.  מאמר די מעניין שנותן הסבר מסקרן למה ואיך ICL עובד. יהיה מעניין לראות מה קורה במקרים שמודל שפה נכשל ב-ICL אם מופעל בצורה הרגילה. האם ההפרדה הזו תישמר?'
$d.Paragraphs.Item(17).Range.InsertParagraphAfter()

Write-Output ("Final paragraph count=" + $d.Paragraphs.Count)
